$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "61.940.05"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "3.415.18"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.89%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.414.53"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +9.99%  "
$ws.Range("E12").Value = "  +6.75%  "
$ws.Range("D13").Value = "4.002.71"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("E15").Value = "  +8.59%  "
$ws.Range("D16").Value = "3.415.84"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("E17").Value = "  +6.17%  "
$ws.Range("D18").Value = "62.018.30"
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.41%  "
$ws.Range("E20").Value = "  +5.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.00%  "
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("D24").Value = "3.553.98"
$ws.Range("E24").Value = "  +3.91%  "
$ws.Range("E25").Value = "  +18.13%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.57%  "
$ws.Range("E28").Value = "  +6.54%  "
$ws.Range("E29").Value = "  +9.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("E32").Value = "  +5.84%  "
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("D34").Value = "3.447.51"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("E39").Value = "  +6.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  +16.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.791"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.15%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("E46").Value = "  +5.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.32%  "
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.99%  "
$ws.Range("D51").Value = "2.381.53"
$ws.Range("E51").Value = "  +11.50%  "
